$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values 5 in C10, D10, E10 (row for Губанов Арсений)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 5

# Set values 5 in C18, D18, E18 (row for Рублева Маргарита)
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 5

# Update the active cell/selection to E18
$ws.Range("E18").Select()
